# Widen column B (the "name" column) and tighten the data-row heights,
# matching the re-formatting done after migrating the data export from
# pip to uv (longer/differently wrapped name values, slightly shorter
# rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B needs to grow from ~13.58 chars to ~31.72 chars so the longer
# "name" values fit without wrapping.
$ws.Columns.Item(2).ColumnWidth = 30.83

# Rows 2-35 (every data row, i.e. everything below the header) shrink
# from 18.75pt to 17.25pt.
$ws.Rows("2:35").RowHeight = 17.25
